# Append new daily rows (329-343) to Sheet1, continuing the time series
# through 2021-08-09 (commit: "aggiornamento fino a 9 agosto 2021").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data scraped from the diff: row, date-serial, col B, col C, col D
$newRows = @(
    @(329, 44403, 0, 4, 95.30617107457708),
    @(330, 44404, 0, 4, 95.30617107457708),
    @(331, 44405, 0, 4, 95.30617107457708),
    @(332, 44406, 0, 4, 95.30617107457708),
    @(333, 44407, 0, 3, 71.47962830593281),
    @(334, 44408, 0, 1, 23.82654276864427),
    @(335, 44409, 0, 0, 0),
    @(336, 44410, 2, 2, 47.65308553728854),
    @(337, 44411, 0, 2, 47.65308553728854),
    @(338, 44412, 0, 2, 47.65308553728854),
    @(339, 44413, 1, 3, 71.47962830593281),
    @(340, 44414, 1, 4, 95.30617107457708),
    @(341, 44415, 0, 4, 95.30617107457708),
    @(342, 44416, 1, 5, 119.1327138432213),
    @(343, 44417, 0, 3, 71.47962830593281)
)

$xlPasteFormats = -4122

foreach ($r in $newRows) {
    $row = $r[0]

    # Column A keeps the same date format / alignment / border as the row above it.
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats) | Out-Null

    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
